$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: single odds update
$ws.Range("AC12").Value = 21

# Row 15: odds refresh (bookmaker re-priced the match)
$ws.Range("G15").Value  = 1.67
$ws.Range("H15").Value  = 3.6
$ws.Range("I15").Value  = 5.25
$ws.Range("J15").Value  = 2.38
$ws.Range("K15").Value  = 2
$ws.Range("L15").Value  = 6
$ws.Range("M15").Value  = 1.08
$ws.Range("N15").Value  = 8
$ws.Range("Q15").Value  = 2.4
$ws.Range("R15").Value  = 1.53
$ws.Range("S15").Value  = 5
$ws.Range("T15").Value  = 1.17
$ws.Range("U15").Value  = 1.53
$ws.Range("V15").Value  = 2.38
$ws.Range("W15").Value  = 2.38
$ws.Range("X15").Value  = 1.53
$ws.Range("Z15").Value  = 6.5
$ws.Range("AA15").Value = 9
$ws.Range("AB15").Value = 12
$ws.Range("AC15").Value = 17
$ws.Range("AF15").Value = 7.5
$ws.Range("AG15").Value = 23
$ws.Range("AH15").Value = 81
$ws.Range("AI15").Value = 10
$ws.Range("AJ15").Value = 26
$ws.Range("AK15").Value = 19
$ws.Range("AL15").Value = 51
$ws.Range("AN15").Value = 51
$ws.Range("AP15").Value = 1.85
$ws.Range("AQ15").Value = 2
$ws.Range("AR15").Value = 3.85
$ws.Range("AS15").Value = 1.25
